# Add motorbikes freight data: row 7 ("motorbikes") on the freight sheet
# was all zeros; it now mirrors the HDVs row (row 3) via formulas.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Make the freight sheet the active tab/sheet (matches the workbook's
# new activeTab / tabSelected state).
$ws.Activate()

# B7 was entered as a single formula...
$ws.Range("B7").Formula = "=B3"
# ...then C7 was entered and filled right through H7, producing one
# shared formula group covering C7:H7.
$ws.Range("C7:H7").Formula = "=C3"

# Leave the selection on the newly-filled row.
$ws.Range("B7:H7").Select()
